$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 190, shifting rows 190:263 down to 191:264
$ws.Rows.Item(190).Insert()

# The cells that just moved from old row 190 to new row 191 still carry the
# correct data for every column except D, K, L, M, P (which get brand new
# values). Copy those unchanged columns back up into the freshly inserted
# row 190, cell by cell (using Value2 to read the underlying value).
$cols = @(1,2,3,5,6,7,8,9,10,14,15,17,18)
foreach ($col in $cols) {
    $ws.Cells.Item(190, $col).Value = $ws.Cells.Item(191, $col).Value2
}

# New values for this inserted record
$ws.Cells.Item(190, 4).Value = 45009
$ws.Cells.Item(190, 4).NumberFormat = $ws.Cells.Item(191, 4).NumberFormat
$ws.Cells.Item(190, 11).Value = 15000
$ws.Cells.Item(190, 12).Value = 16000
$ws.Cells.Item(190, 13).Value = 15500
$ws.Cells.Item(190, 16).Value = 1550
